$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 0.89812588
$ws.Range("D3").Value = 1.9753469
